# Applies the roster re-shuffle described in the commit diff:
# rows are renumbered/reshuffled so several players now sit on different
# rows, and some players change team/position (reflecting trades).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("Darius Garland", "PG", "Cleveland Cavaliers"),
  @("Klay Thompson", "SG,SF", "Dallas Mavericks"),
  @("Stephen Curry", "PG,SG", "Golden State Warriors"),
  @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
  @("OG Anunoby", "SF,PF", "New York Knicks"),
  @("Keegan Murray", "SF,PF", "Sacramento Kings"),
  @("Kevin Durant", "SF,PF", "Phoenix Suns"),
  @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
  @("Jarrett Allen", "C", "Cleveland Cavaliers"),
  @("Mark Williams", "C", "Charlotte Hornets"),
  @("Trey Murphy III", "SG,SF,PF", "New Orleans Pelicans"),
  @("Jalen Duren", "C", "Detroit Pistons"),
  @("Franz Wagner", "SF,PF", "Orlando Magic"),
  @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
  @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
  @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
